$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 50")

# Row 7: center E7, center+wrap F7 (formatting update to match the other weeks)
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").VerticalAlignment = -4108
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("F7").WrapText = $true

# Row 8: new activity entry - begin/end time, description text
$ws.Range("C8").NumberFormat = "h:mm;@"
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("C8").VerticalAlignment = -4108
$ws.Range("C8").Value = 0.38541666666666669

$ws.Range("D8").NumberFormat = "h:mm;@"
$ws.Range("D8").HorizontalAlignment = -4108
$ws.Range("D8").VerticalAlignment = -4108
$ws.Range("D8").Value = 0.4236111111111111

$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").VerticalAlignment = -4108

$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("F8").WrapText = $true
$ws.Range("F8").Value = "StartScene achtergrond toegevoegd en begin gemaakt aan title"

$ws.Rows.Item(8).RowHeight = 28.5

# Rows 9-11: center E, center+wrap F
$ws.Range("E9:E11").HorizontalAlignment = -4108
$ws.Range("E9:E11").VerticalAlignment = -4108
$ws.Range("F9:F17").HorizontalAlignment = -4108
$ws.Range("F9:F17").WrapText = $true

# Page setup to match the rest of the workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to F18, like a user who just finished row 8 and tabbed to the bottom
$ws.Range("F18").Select()
